# Commit: removed "error" in excel
# Every cell that displayed the text "error" (State_new column, col C) on the
# three worksheets gets its text shortened to "e". This causes the old
# "error" shared-string slot to become unused (and get dropped), with the
# following strings shifting down and "e" appended at the end - which is
# exactly what the target workbook's sharedStrings.xml shows.
#
# The commit also left the workbook with a different sheet/cell selection
# (the author apparently clicked around after editing): Multiplication
# becomes the active/selected tab instead of Invert, and each sheet's
# remembered selection cell changes too.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Addition")
$ws1.Range("C3").Value = "e"
$ws1.Range("C4").Value = "e"
$ws1.Range("C11").Value = "e"
$ws1.Range("C14").Value = "e"
$ws1.Range("C15").Value = "e"
$ws1.Range("C16").Value = "e"

$ws2 = $wb.Worksheets.Item("Invert")
$ws2.Range("C3").Value = "e"
$ws2.Range("C4").Value = "e"
$ws2.Range("C8").Value = "e"

$ws3 = $wb.Worksheets.Item("Multiplication")
$ws3.Range("C3").Value = "e"
$ws3.Range("C4").Value = "e"
$ws3.Range("C11").Value = "e"

# Restore per-sheet selections, ending on Multiplication so it becomes the
# active tab (matches activeTab going from 1 -> 2 and tabSelected moving
# from Invert to Multiplication).
$ws1.Range("D19").Select()
$ws2.Range("D15").Select()
$ws3.Range("C18").Select()
